$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.373.19'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '1.622.59'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '212.12'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.485'
$ws.Range("E7").Value = '  +1.22%  '
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").Value = '18.86'
$ws.Range("E10").Value = '  +3.84%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '1.847.94'
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("D13").Value = '1.634.57'
$ws.Range("E13").Value = '  +2.70%  '
$ws.Range("D14").Value = '4.01'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").Value = '0.518'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '26.364.21'
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").Value = '62.48'
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = '202.63'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").Value = '4.27'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = '9.33'
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").Value = '6.04'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '1.88'
$ws.Range("E24").Value = '  -3.20%  '
$ws.Range("D25").Value = '144.47'
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").Value = '15.20'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").Value = '0.0523'
$ws.Range("E30").Value = '  +10.03%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("D34").Value = '1.49'
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("E35").Value = '  +2.08%  '
$ws.Range("D36").Value = '1.177.20'
$ws.Range("E36").Value = '  +4.34%  '
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").Value = '0.807'
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").Value = '0.498'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("E42").Value = '  +4.70%  '
$ws.Range("D43").Value = '0.785'
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("D44").Value = '1.759.78'
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").Value = '92.74'
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '1.52'
$ws.Range("E46").Value = '  +2.45%  '
$ws.Range("D47").Value = '53.86'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = '0.410'
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").Value = '7.30'
$ws.Range("E51").Value = '  +1.49%'
